$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the last data row (was row 5: Sending=MuSCs, Target=Resolving-Mac) ---
# The whole row is gone in the new TPM export, remaining rows shift up.
$ws.Rows.Item(5).Delete()

# --- Sending cluster column (A) for the remaining 3 data rows: MuSCs -> Resolving-Mac ---
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("A4").Value = "Resolving-Mac"

# --- Row 2 (Target cluster: ECs) updated TPM-derived values ---
$ws.Range("G2").Value = 0.02530666666666667
$ws.Range("H2").Value = 0.07592
$ws.Range("M2").Value = 0.8155003333333334
$ws.Range("N2").Value = 2.446501
$ws.Range("O2").Value = 0.1910612426590028
$ws.Range("P2").Value = 0.1910612426590029
$ws.Range("Q2").Value = 0.02063759510222223
$ws.Range("R2").Value = 0.18573835592
$ws.Range("S2").Value = 0.1910612426590028
$ws.Range("T2").Value = 0.1910612426590029

# --- Row 3 (Target cluster: FAPs) updated TPM-derived values ---
$ws.Range("G3").Value = 0.02530666666666667
$ws.Range("H3").Value = 0.07592
$ws.Range("O3").Value = 0.7809105179307759
$ws.Range("P3").Value = 0.780910517930776
$ws.Range("Q3").Value = 0.0843505195288889
$ws.Range("R3").Value = 0.7591546757600001
$ws.Range("S3").Value = 0.7809105179307759
$ws.Range("T3").Value = 0.780910517930776

# --- Row 4 (Target cluster: MuSCs) updated TPM-derived values ---
$ws.Range("G4").Value = 0.02530666666666667
$ws.Range("H4").Value = 0.07592
$ws.Range("M4").Value = 0.119632
$ws.Range("N4").Value = 0.358896
$ws.Range("O4").Value = 0.02802823941022116
$ws.Range("P4").Value = 0.02802823941022117
$ws.Range("Q4").Value = 0.003027487146666667
$ws.Range("R4").Value = 0.02724738432
$ws.Range("S4").Value = 0.02802823941022116
$ws.Range("T4").Value = 0.02802823941022117
